# Thong ke hop tac: them so lieu vao cuoi hai gach dau dong.
$d = $word.ActiveDocument

function Append-ItalicText($searchText, $appendText) {
    $rng = $d.Content.Duplicate
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Khong tim thay doan van ban: $searchText"
    }
    $rng.Collapse(0)  # wdCollapseEnd
    $rng.InsertAfter($appendText)
    $rng.Font.Italic = $true
}

Append-ItalicText "Số lượng collaborator tham dự (của dự án gốc)" ": 3"
Append-ItalicText "Số lượt commit" ": 15"
